$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new blank column before column N ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$mWidth = $wsRepay.Columns("M:M").ColumnWidth
$wsRepay.Columns("N:N").Insert()
$wsRepay.Columns("N:N").ColumnWidth = $mWidth

# --- Update selections / active sheet / active tab ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F5").Select()

$wsRepay.Range("L15").Select()
$wsRepay.Activate()
